# Scheduled-runner update: refresh Sheets price/profit columns (H-N) with
# latest market data snapshot across ALC/ARM/CRP/CUL/GSM/LTW/WVR tabs.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4708.6665
$ws.Range("I113").Value = 4317.3335
$ws.Range("K113").Value = 4317.3335
$ws.Range("M113").Value = -1063.3335
$ws.Range("H138").Value = 2793.8
$ws.Range("I138").Value = 2085.65
$ws.Range("K138").Value = 6256.950000000001
$ws.Range("M138").Value = -1116.950000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4107.6924
$ws.Range("I61").Value = 2980
$ws.Range("J61").Value = 12753.333
$ws.Range("K61").Value = 2980
$ws.Range("L61").Value = 12753.333
$ws.Range("M61").Value = -2768
$ws.Range("N61").Value = -13177.333
$ws.Range("H74").Value = 4456.643
$ws.Range("I74").Value = 1174.125
$ws.Range("J74").Value = 8833.333000000001
$ws.Range("K74").Value = 1174.125
$ws.Range("L74").Value = 8833.333000000001
$ws.Range("M74").Value = -300.125
$ws.Range("N74").Value = -10581.333
$ws.Range("H77").Value = 4456.643
$ws.Range("I77").Value = 1174.125
$ws.Range("J77").Value = 8833.333000000001
$ws.Range("K77").Value = 5870.625
$ws.Range("L77").Value = 44166.665
$ws.Range("M77").Value = -1502.625
$ws.Range("N77").Value = -52902.665
$ws.Range("H132").Value = 7933.391
$ws.Range("I132").Value = 4185.875
$ws.Range("K132").Value = 12557.625
$ws.Range("M132").Value = -10027.625
$ws.Range("H134").Value = 120000
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("H136").Value = 4107.6924
$ws.Range("I136").Value = 2980
$ws.Range("J136").Value = 12753.333
$ws.Range("K136").Value = 8940
$ws.Range("L136").Value = 38259.999
$ws.Range("M136").Value = -6390
$ws.Range("N136").Value = -43359.999
$ws.Range("H139").Value = 89998.60000000001
$ws.Range("J139").Value = 89998.60000000001
$ws.Range("L139").Value = 89998.60000000001
$ws.Range("N139").Value = -100278.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3974.4443
$ws.Range("I58").Value = 2622.48
$ws.Range("K58").Value = 2622.48
$ws.Range("M58").Value = -2419.48
$ws.Range("H103").Value = 7584
$ws.Range("I103").Value = 7584
$ws.Range("K103").Value = 7584
$ws.Range("M103").Value = -6412
$ws.Range("H136").Value = 3974.4443
$ws.Range("I136").Value = 2622.48
$ws.Range("K136").Value = 7867.440000000001
$ws.Range("M136").Value = -5317.440000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 346.625
$ws.Range("I5").Value = 223
$ws.Range("J5").Value = 408.4375
$ws.Range("K5").Value = 669
$ws.Range("L5").Value = 1225.3125
$ws.Range("M5").Value = -557
$ws.Range("N5").Value = -1449.3125
$ws.Range("H57").Value = 1697.8
$ws.Range("I57").Value = 694.5
$ws.Range("K57").Value = 2083.5
$ws.Range("M57").Value = -1524.5
$ws.Range("H58").Value = 1467.5
$ws.Range("J58").Value = 2003
$ws.Range("L58").Value = 6009
$ws.Range("N58").Value = -6265
$ws.Range("H114").Value = 627.7692
$ws.Range("J114").Value = 999.5
$ws.Range("L114").Value = 2998.5
$ws.Range("N114").Value = -9506.5
$ws.Range("H117").Value = 3218.2856
$ws.Range("I117").Value = 3257
$ws.Range("J117").Value = 3166.6667
$ws.Range("K117").Value = 9771
$ws.Range("L117").Value = 9500.000100000001
$ws.Range("M117").Value = -6329
$ws.Range("N117").Value = -16384.0001
$ws.Range("H122").Value = 14192.143
$ws.Range("J122").Value = 4983
$ws.Range("L122").Value = 44847
$ws.Range("N122").Value = -49747
$ws.Range("H135").Value = 346.625
$ws.Range("I135").Value = 223
$ws.Range("J135").Value = 408.4375
$ws.Range("K135").Value = 2007
$ws.Range("L135").Value = 3675.9375
$ws.Range("M135").Value = 528
$ws.Range("N135").Value = -8745.9375
$ws.Range("H139").Value = 5034.5386
$ws.Range("I139").Value = 2686.625
$ws.Range("K139").Value = 8059.875
$ws.Range("M139").Value = -2919.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 53027.79
$ws.Range("J2").Value = 111155.78
$ws.Range("L2").Value = 111155.78
$ws.Range("N2").Value = -111381.78
$ws.Range("H9").Value = 2002.6666
$ws.Range("I9").Value = 1504
$ws.Range("K9").Value = 1504
$ws.Range("M9").Value = -1334
$ws.Range("H113").Value = 37776
$ws.Range("I113").Value = 4797.6
$ws.Range("J113").Value = 78999
$ws.Range("K113").Value = 4797.6
$ws.Range("L113").Value = 78999
$ws.Range("M113").Value = -2627.6
$ws.Range("N113").Value = -83339
$ws.Range("H132").Value = 8360
$ws.Range("I132").Value = 5826.231
$ws.Range("K132").Value = 17478.693
$ws.Range("M132").Value = -14948.693

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2622.5
$ws.Range("I68").Value = 2711.4285
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 2711.4285
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -1962.4285
$ws.Range("N68").Value = -3498
$ws.Range("H71").Value = 2622.5
$ws.Range("I71").Value = 2711.4285
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 13557.1425
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -9813.1425
$ws.Range("N71").Value = -17488
$ws.Range("H132").Value = 6139.2144
$ws.Range("I132").Value = 4828
$ws.Range("J132").Value = 7450.4287
$ws.Range("K132").Value = 14484
$ws.Range("L132").Value = 22351.2861
$ws.Range("M132").Value = -11954
$ws.Range("N132").Value = -27411.2861
$ws.Range("H138").Value = 25000
$ws.Range("J138").Value = 25000
$ws.Range("L138").Value = 25000
$ws.Range("N138").Value = -35280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 8266.786
$ws.Range("I96").Value = 1019.7143
$ws.Range("J96").Value = 15513.857
$ws.Range("K96").Value = 1019.7143
$ws.Range("L96").Value = 15513.857
$ws.Range("M96").Value = 353.2857
$ws.Range("N96").Value = -18259.857
$ws.Range("H100").Value = 288.7143
$ws.Range("I100").Value = 249.11111
$ws.Range("J100").Value = 360
$ws.Range("K100").Value = 498.22222
$ws.Range("L100").Value = 720
$ws.Range("M100").Value = 42.77778000000001
$ws.Range("N100").Value = -1802
$ws.Range("H126").Value = 3589.6924
$ws.Range("I126").Value = 3589.6924
$ws.Range("K126").Value = 10769.0772
$ws.Range("M126").Value = -8299.0772
$ws.Range("H132").Value = 4561.7964
$ws.Range("I132").Value = 2825.9736
$ws.Range("J132").Value = 7702.8096
$ws.Range("K132").Value = 8477.9208
$ws.Range("L132").Value = 23108.4288
$ws.Range("M132").Value = -5947.9208
$ws.Range("N132").Value = -28168.4288
